# The document has a single comment ("Adam Fyne" / "Comment on 'Word'.")
# whose range anchors the word "with" ("Word [with] a comment."). The
# target state of the file turns this into an *empty annotation mark*:
# the commented word is removed so the comment's range collapses to a
# single point (commentRangeStart/commentRangeEnd back-to-back) while the
# comment itself, its reference and all other document content/styles
# stay untouched.
$d = $word.ActiveDocument

# Locate the (only) comment and delete the text it is anchored to -
# this is exactly what happens when a user selects the annotated text
# for a comment and presses Delete: the comment survives, but now
# anchors an empty range ("empty annotation mark").
For ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $c = $d.Comments.Item($i)
    $c.Scope.Delete()
}
